$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (STATLED1): Kingbright APHHS1005QBC/D (Blue 465nm) -> Inolux IN-S42BT5B (Blue 467nm)
$ws.Range("B12").Value = "IN-S42BT5B"
$ws.Range("A12").Value = "Inolux"
$ws.Range("E12").Value = "Blue 467nm LED Indication - Discrete 2.9V 0402 (1005 Metric)"

# Row 13 (STATLED2): Kingbright APHHS1005SECK (Orange 601nm) -> Inolux IN-S42BT5A (Amber 605nm)
$ws.Range("B13").Value = "IN-S42BT5A"
$ws.Range("A13").Value = "Inolux"
$ws.Range("E13").Value = "Amber 605nm LED Indication - Discrete 2V 0402 (1005 Metric)"

# Restore the view/selection state recorded in the saved workbook
$ws.Range("E13").Select()
